$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Find-ParaIndex($text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text.StartsWith($text)) {
            return $i
        }
    }
    return -1
}

# Replace the paragraph whose text starts with $text with the supplied
# paragraph-level OOXML (one or more <w:p>...</w:p> fragments). Returns the
# index of the (first) newly inserted paragraph.
function Replace-Paragraph($text, $xml) {
    $idx = Find-ParaIndex($text)
    $old = $d.Paragraphs.Item($idx)
    $old.Range.InsertParagraphAfter()
    $slot = $d.Paragraphs.Item($idx + 1)
    $slot.Range.InsertXML($xml)
    $d.Paragraphs.Item($idx).Range.Delete()
    return $idx
}

# Insert new paragraph(s) (OOXML) right after the paragraph whose text
# starts with $text.
function Insert-ParagraphsAfter($text, $xml) {
    $idx = Find-ParaIndex($text)
    $anchor = $d.Paragraphs.Item($idx)
    $anchor.Range.InsertParagraphAfter()
    $slot = $d.Paragraphs.Item($idx + 1)
    $slot.Range.InsertXML($xml)
}

# Insert new paragraph(s) (OOXML) right before the paragraph whose text
# starts with $text.
function Insert-ParagraphsBefore($text, $xml) {
    $idx = Find-ParaIndex($text)
    $anchor = $d.Paragraphs.Item($idx)
    $anchor.Range.InsertParagraphBefore()
    $slot = $d.Paragraphs.Item($idx)
    $slot.Range.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 1) Concatenation paragraph: remove the stray _GoBack bookmark that used to
#    sit between "ings or several" and " lists".
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks
if ($bm.Exists("_GoBack")) {
    $bm.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) Compiler paragraph: split the single run
#      "are system which converts source code to lower-level"
#    into two runs:
#      "are system which converts"
#      " code to lower-level"
# ---------------------------------------------------------------------------
$compilerXml = '<w:p ' + $wNs + '>' + `
    '<w:r><w:rPr><w:b/></w:rPr><w:t>Compiler</w:t></w:r>' + `
    '<w:r><w:t>: a softw</w:t></w:r>' + `
    '<w:r><w:t>are system which converts</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> code to lower-level</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> code</w:t></w:r>' + `
    '</w:p>'
Replace-Paragraph "Compiler:" $compilerXml | Out-Null

# ---------------------------------------------------------------------------
# 3) Insert "Function" and "Graphical User Interface" paragraphs right after
#    the Documentation paragraph (i.e. right before "Integrated Development
#    Environment (IDE)").
# ---------------------------------------------------------------------------
$funcGuiXml = '<w:p ' + $wNs + '>' + `
    '<w:r><w:rPr><w:b/></w:rPr><w:t>Function</w:t></w:r>' + `
    '<w:r><w:t>: a set of instructions that can be referenced by a name</w:t></w:r>' + `
    '</w:p>' + `
    '<w:p ' + $wNs + '>' + `
    '<w:r><w:rPr><w:b/></w:rPr><w:t>Graphical User Interface</w:t></w:r>' + `
    '<w:r><w:t>: an interface which allows the use of electronic devices via graphical icons and visual cues</w:t></w:r>' + `
    '</w:p>'
Insert-ParagraphsAfter "Documentation:" $funcGuiXml

# ---------------------------------------------------------------------------
# 4) Insert "Interpreter", "Library" and "Loop" paragraphs right after the
#    Integrated Development Environment (IDE) paragraph.
# ---------------------------------------------------------------------------
$interpLibLoopXml = '<w:p ' + $wNs + '>' + `
    '<w:r><w:rPr><w:b/></w:rPr><w:t>Interpreter</w:t></w:r>' + `
    '<w:r><w:t>: a software system which converts code to lower-level code on-the-fly</w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
    '</w:p>' + `
    '<w:p ' + $wNs + '>' + `
    '<w:r><w:rPr><w:b/></w:rPr><w:t>Library</w:t></w:r>' + `
    '<w:r><w:t>: a collection of implementations for the purpose of general reuse</w:t></w:r>' + `
    '</w:p>' + `
    '<w:p ' + $wNs + '>' + `
    '<w:r><w:rPr><w:b/></w:rPr><w:t>Loop</w:t></w:r>' + `
    '<w:r><w:t>: a repeated sequence of instructions until some condition is met</w:t></w:r>' + `
    '</w:p>'
Insert-ParagraphsAfter "Integrated Development Environment (IDE):" $interpLibLoopXml

# ---------------------------------------------------------------------------
# 5) Insert a new "Method" paragraph (": a function that is associated with
#    an object") right before the existing "Method Overloading" paragraph.
# ---------------------------------------------------------------------------
$newMethodXml = '<w:p ' + $wNs + '>' + `
    '<w:r><w:rPr><w:b/></w:rPr><w:t>Method</w:t></w:r>' + `
    '<w:r><w:t>: a function that is associated with an object</w:t></w:r>' + `
    '</w:p>'
Insert-ParagraphsBefore "Method Overloading:" $newMethodXml

# ---------------------------------------------------------------------------
# 6) Remove the old "Library" and "Loop" paragraphs that used to sit right
#    after "Method Overriding".
# ---------------------------------------------------------------------------
$moIdx = Find-ParaIndex("Method Overriding:")
$libPara = $d.Paragraphs.Item($moIdx + 1)
$loopPara = $d.Paragraphs.Item($moIdx + 2)
$delRange = $d.Range($libPara.Range.Start, $loopPara.Range.End)
$delRange.Delete()

Write-Output "done"
